# Remove the "kelas_id" column (A) and the "angkatan" column (H, which
# becomes G after the first deletion). Remaining columns shift left so the
# sheet goes from A1:H2 to A1:F2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Delete()
$ws.Columns("G").Delete()

# Match the author's saved selection state (A1:A2).
$ws.Range("A1:A2").Select()
